$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("I3").Value = 0.1
$ws3.Range("I3").NumberFormat = "0%"
$ws3.Range("I3").Style = "Percent"
